$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference colors already used on the sheet for the three status values
# (TODO = red, DONE = dark green, SATISFACTORY = light green)
$todoColor = $ws.Range("H4").Interior.Color()
$doneColor = $ws.Range("H8").Interior.Color()
$satColor  = $ws.Range("H7").Interior.Color()

# Row 9: "Allow user input" task now links to the maze-resize note and is
# reclassified from DONE to SATISFACTORY
$ws.Range("L9").Value = "Allow user to resize the dimensions of the maze"
$ws.Range("M9").Value = "SATISFACTORY"
$ws.Range("M9").Interior.Color = $satColor

# Rows 12 & 13: both move from TODO to DONE
$ws.Range("M12").Value = "DONE"
$ws.Range("M12").Interior.Color = $doneColor

$ws.Range("M13").Value = "DONE"
$ws.Range("M13").Interior.Color = $doneColor

# Row 14: new entry in the Networks table
$ws.Range("K14").Value = "User move start/ end position"
$ws.Range("M14").Value = "TODO"
$ws.Range("M14").Interior.Color = $todoColor

# Row 15: new entry in the Networks table
$ws.Range("K15").Value = "Server console output"
$ws.Range("L15").Value = "For each different packet"
$ws.Range("M15").Value = "TODO"
$ws.Range("M15").Interior.Color = $todoColor

# Update the current selection to match the author's saved view
$ws.Range("F32").Select() | Out-Null
